# chapter 9 auth and server setup sessions.docx
# Commit: "production mode dist dir"
#
# 1) Collapse several runs (that were split apart only because of
#    w:proofErr spell/gram markers) back into single runs with the
#    same visible text.
# 2) Append a short new section after the JSON cookie example that
#    explains how to inspect the build output tree (`tree -F ./` /
#    `tree -F ./dist`).

$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Part 1: merge split runs into single runs (same text, no meaning
# change - just removes the spellcheck run-splitting).
# ---------------------------------------------------------------

# "Cat sessions/*.json | jq " (paragraph 2 - leaves the rest of the
# paragraph, including the curly-quote '...' run and gramStart/gramEnd,
# untouched).
$m1 = $d.Range(28, 53)
$m1.Text = "Cat sessions/*.json | jq "

# "originalMaxAge": null,
$m2 = $d.Range(102, 125)
$m2.Text = """originalMaxAge"": null,"

# "httpOnly": true,
$m3 = $d.Range(143, 160)
$m3.Text = """httpOnly"": true,"

# "path": "/"
$m4 = $d.Range(161, 172)
$m4.Text = """path"": ""/"""

# "         ""__lastAccess"": 1499592937433"
$m5 = $d.Range(176, 214)
$m5.Text = "         ""__lastAccess"": 1499592937433"

# ---------------------------------------------------------------
# Part 2: append the new "tree" section after the closing "}" of the
# session JSON block (still paragraph 11 - none of the merges above
# changed paragraph count or the visible text length of paragraph 11).
# ---------------------------------------------------------------

$lastPara = $d.Paragraphs.Item(11)
$lastPara.Range.InsertParagraphAfter()
$insertionPoint = $d.Paragraphs.Item(12).Range

$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="DejaVuSansMono" w:hAnsi="DejaVuSansMono" w:cs="DejaVuSansMono"/><w:kern w:val="0"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="DejaVuSansMono" w:hAnsi="DejaVuSansMono" w:cs="DejaVuSansMono"/><w:kern w:val="0"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="DejaVuSansMono" w:hAnsi="DejaVuSansMono" w:cs="DejaVuSansMono"/><w:kern w:val="0"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>We can check the t</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="DejaVuSansMono" w:hAnsi="DejaVuSansMono" w:cs="DejaVuSansMono"/><w:kern w:val="0"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>ree structure in the files like this linux:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="DejaVuSansMono" w:hAnsi="DejaVuSansMono" w:cs="DejaVuSansMono"/><w:kern w:val="0"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="DejaVuSansMono" w:hAnsi="DejaVuSansMono" w:cs="DejaVuSansMono"/><w:kern w:val="0"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>Tree -</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="DejaVuSansMono" w:hAnsi="DejaVuSansMono" w:cs="DejaVuSansMono"/><w:kern w:val="0"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>F .</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="DejaVuSansMono" w:hAnsi="DejaVuSansMono" w:cs="DejaVuSansMono"/><w:kern w:val="0"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>/ ---- follow by the file we want to see if it exist</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="DejaVuSansMono" w:hAnsi="DejaVuSansMono" w:cs="DejaVuSansMono"/><w:kern w:val="0"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">  Tree -</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="DejaVuSansMono" w:hAnsi="DejaVuSansMono" w:cs="DejaVuSansMono"/><w:kern w:val="0"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>F .</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="DejaVuSansMono" w:hAnsi="DejaVuSansMono" w:cs="DejaVuSansMono"/><w:kern w:val="0"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>/dist</w:t></w:r></w:p>
'@

$insertionPoint.InsertXML($newXml)

Write-Output "done"
